$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the brand list: old 4-row table (BrandName/Select header + 3 data
# rows) becomes a 2-row table with just a header cell and one brand.
$ws.Range("A1").Value = "Boroline"
$ws.Range("B1").Value = $null

$ws.Range("A2").Value = "Dettol"
$ws.Range("A2").Font.Bold = $true

$ws.Range("B2").Value = $null
$ws.Range("A3:B4").Value = $null

$ws.Range("A2").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
